$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.71
$ws.Range("C2").Value = 0.6899999999999999
$ws.Range("D2").Value = 0.8
$ws.Range("E2").Value = 0.95
$ws.Range("F2").Value = 0.98
$ws.Range("G2").Value = 0.8100000000000001
$ws.Range("H2").Value = 0.93
$ws.Range("I2").Value = 0.96
$ws.Range("J2").Value = 0.9
$ws.Range("K2").Value = 0.86

$ws.Range("B3").Value = 0.76
$ws.Range("C3").Value = 0.68
$ws.Range("D3").Value = 0.75
$ws.Range("E3").Value = 0.91
$ws.Range("F3").Value = 0.99
$ws.Range("G3").Value = 0.86
$ws.Range("H3").Value = 0.92
$ws.Range("I3").Value = 0.97
$ws.Range("J3").Value = 0.84
$ws.Range("K3").Value = 0.85

$ws.Range("B4").Value = 0.78
$ws.Range("C4").Value = 0.65
$ws.Range("D4").Value = 0.8100000000000001
$ws.Range("E4").Value = 0.86
$ws.Range("F4").Value = 0.98
$ws.Range("G4").Value = 0.85
$ws.Range("H4").Value = 0.91
$ws.Range("I4").Value = 0.95
$ws.Range("J4").Value = 0.88
$ws.Range("K4").Value = 0.85

$ws.Range("B5").Value = 0.7
$ws.Range("C5").Value = 0.7
$ws.Range("D5").Value = 0.79
$ws.Range("E5").Value = 0.93
$ws.Range("F5").Value = 0.98
$ws.Range("G5").Value = 0.86
$ws.Range("H5").Value = 0.9
$ws.Range("I5").Value = 0.95
$ws.Range("J5").Value = 0.84
$ws.Range("K5").Value = 0.85

$ws.Range("B6").Value = 0.74
$ws.Range("C6").Value = 0.65
$ws.Range("D6").Value = 0.68
$ws.Range("E6").Value = 0.91
$ws.Range("F6").Value = 0.98
$ws.Range("G6").Value = 0.86
$ws.Range("H6").Value = 0.92
$ws.Range("I6").Value = 0.97
$ws.Range("J6").Value = 0.87
$ws.Range("K6").Value = 0.84

$ws.Range("B7").Value = 0.78
$ws.Range("C7").Value = 0.71
$ws.Range("D7").Value = 0.73
$ws.Range("E7").Value = 0.9399999999999999
$ws.Range("F7").Value = 0.97
$ws.Range("G7").Value = 0.88
$ws.Range("H7").Value = 0.9399999999999999
$ws.Range("I7").Value = 0.95
$ws.Range("J7").Value = 0.88
$ws.Range("K7").Value = 0.86

$ws.Range("B8").Value = 0.76
$ws.Range("C8").Value = 0.71
$ws.Range("D8").Value = 0.82
$ws.Range("E8").Value = 0.87
$ws.Range("F8").Value = 0.98
$ws.Range("G8").Value = 0.88
$ws.Range("H8").Value = 0.93
$ws.Range("I8").Value = 0.96
$ws.Range("J8").Value = 0.86
$ws.Range("K8").Value = 0.86

$ws.Range("B9").Value = 0.74
$ws.Range("C9").Value = 0.71
$ws.Range("D9").Value = 0.76
$ws.Range("E9").Value = 0.9
$ws.Range("F9").Value = 0.98
$ws.Range("G9").Value = 0.86
$ws.Range("H9").Value = 0.93
$ws.Range("I9").Value = 0.97
$ws.Range("J9").Value = 0.85
$ws.Range("K9").Value = 0.86

$ws.Range("B10").Value = 0.77
$ws.Range("C10").Value = 0.66
$ws.Range("D10").Value = 0.82
$ws.Range("E10").Value = 0.96
$ws.Range("F10").Value = 0.99
$ws.Range("G10").Value = 0.85
$ws.Range("H10").Value = 0.93
$ws.Range("I10").Value = 0.93
$ws.Range("J10").Value = 0.89
$ws.Range("K10").Value = 0.87

$ws.Range("B11").Value = 0.72
$ws.Range("C11").Value = 0.68
$ws.Range("D11").Value = 0.79
$ws.Range("E11").Value = 0.91
$ws.Range("F11").Value = 0.98
$ws.Range("G11").Value = 0.85
$ws.Range("H11").Value = 0.91
$ws.Range("I11").Value = 0.96
$ws.Range("J11").Value = 0.85
$ws.Range("K11").Value = 0.85

$ws.Range("B12").Value = 0.74
$ws.Range("C12").Value = 0.7
$ws.Range("D12").Value = 0.83
$ws.Range("E12").Value = 0.93
$ws.Range("F12").Value = 0.99
$ws.Range("G12").Value = 0.88
$ws.Range("H12").Value = 0.92
$ws.Range("I12").Value = 0.97
$ws.Range("J12").Value = 0.87
$ws.Range("K12").Value = 0.87

$ws.Range("B13").Value = 0.73
$ws.Range("C13").Value = 0.7
$ws.Range("D13").Value = 0.77
$ws.Range("E13").Value = 0.95
$ws.Range("F13").Value = 0.98
$ws.Range("G13").Value = 0.9
$ws.Range("H13").Value = 0.92
$ws.Range("I13").Value = 0.95
$ws.Range("J13").Value = 0.85
$ws.Range("K13").Value = 0.86

$ws.Range("B14").Value = 0.75
$ws.Range("C14").Value = 0.6899999999999999
$ws.Range("D14").Value = 0.82
$ws.Range("E14").Value = 0.9399999999999999
$ws.Range("F14").Value = 0.97
$ws.Range("G14").Value = 0.86
$ws.Range("H14").Value = 0.92
$ws.Range("I14").Value = 0.97
$ws.Range("J14").Value = 0.87
$ws.Range("K14").Value = 0.87

$ws.Range("B15").Value = 0.8
$ws.Range("C15").Value = 0.6899999999999999
$ws.Range("D15").Value = 0.74
$ws.Range("E15").Value = 0.96
$ws.Range("F15").Value = 0.98
$ws.Range("G15").Value = 0.87
$ws.Range("H15").Value = 0.93
$ws.Range("I15").Value = 0.96
$ws.Range("J15").Value = 0.88
$ws.Range("K15").Value = 0.87

$ws.Range("B16").Value = 0.75
$ws.Range("C16").Value = 0.6899999999999999
$ws.Range("D16").Value = 0.78
$ws.Range("E16").Value = 0.92
$ws.Range("F16").Value = 0.98
$ws.Range("G16").Value = 0.86
$ws.Range("H16").Value = 0.92
$ws.Range("I16").Value = 0.96
$ws.Range("J16").Value = 0.87
$ws.Range("K16").Value = 0.86

